# Update the cryptos price list per the latest GitHub Actions data refresh.

function Set-TextValue {
    param($Cell, $Text)
    # Force the cell to keep a text (string) value even when the text looks numeric,
    # then restore the default "Normal" style so no stray formatting is introduced.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values row by row
Set-TextValue $ws.Cells.Item(2, 4) "29.406.06"
$ws.Cells.Item(2, 5).Value = "  -0.05%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.850.81"
$ws.Cells.Item(3, 5).Value = "  +0.14%  "
Set-TextValue $ws.Cells.Item(4, 4) "0.9998"
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
Set-TextValue $ws.Cells.Item(5, 4) "240.84"
$ws.Cells.Item(5, 5).Value = "  +0.16%  "
$ws.Cells.Item(6, 5).Value = "  -0.53%  "
$ws.Cells.Item(7, 5).Value = "  +0.06%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.07672"
$ws.Cells.Item(8, 5).Value = "  +1.43%  "
$ws.Cells.Item(9, 5).Value = "  -0.70%  "
Set-TextValue $ws.Cells.Item(10, 4) "24.54"
$ws.Cells.Item(10, 5).Value = "  -0.26%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.07747"
$ws.Cells.Item(11, 5).Value = "  +0.67%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.857.39"
$ws.Cells.Item(12, 5).Value = "  +0.40%  "
$ws.Cells.Item(13, 5).Value = "  +9.95%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.021"
$ws.Cells.Item(14, 5).Value = "  +0.60%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.6812"
$ws.Cells.Item(15, 5).Value = "  -0.73%  "
$ws.Cells.Item(16, 5).Value = "  +0.58%  "
Set-TextValue $ws.Cells.Item(17, 4) "2.104.76"
$ws.Cells.Item(17, 5).Value = "  +0.03%  "
Set-TextValue $ws.Cells.Item(18, 4) "6.145"
$ws.Cells.Item(18, 5).Value = "  +0.24%  "
Set-TextValue $ws.Cells.Item(19, 4) "29.454.48"
$ws.Cells.Item(19, 5).Value = "  +0.02%  "
Set-TextValue $ws.Cells.Item(20, 4) "229.25"
$ws.Cells.Item(20, 5).Value = "  +0.21%  "
$ws.Cells.Item(21, 5).Value = "  -0.23%  "
$ws.Cells.Item(22, 5).Value = "  +0.04%  "
Set-TextValue $ws.Cells.Item(23, 4) "7.437"
$ws.Cells.Item(24, 5).Value = "  +0.03%  "
Set-TextValue $ws.Cells.Item(25, 4) "156.94"
$ws.Cells.Item(25, 5).Value = "  -0.12%  "
$ws.Cells.Item(26, 5).Value = "  -0.78%  "
$ws.Cells.Item(27, 5).Value = "  +0.17%  "
Set-TextValue $ws.Cells.Item(28, 4) "17.68"
$ws.Cells.Item(28, 5).Value = "  -0.01%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.317"
$ws.Cells.Item(29, 5).Value = "  +4.04%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.466"
Set-TextValue $ws.Cells.Item(31, 4) "0.05708"
$ws.Cells.Item(31, 5).Value = "  -0.01%  "
Set-TextValue $ws.Cells.Item(32, 4) "4.126"
$ws.Cells.Item(32, 5).Value = "  +0.08%  "
$ws.Cells.Item(33, 5).Value = "  +0.55%  "
$ws.Cells.Item(34, 5).Value = "  +0.16%  "
Set-TextValue $ws.Cells.Item(35, 4) "1.163"
$ws.Cells.Item(35, 5).Value = "  +0.68%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.7041"
$ws.Cells.Item(36, 5).Value = "  -1.49%  "
Set-TextValue $ws.Cells.Item(37, 4) "2.583"
$ws.Cells.Item(37, 5).Value = "  -0.16%  "
Set-TextValue $ws.Cells.Item(38, 4) "2.780"
$ws.Cells.Item(38, 5).Value = "  +0.12%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.01791"
$ws.Cells.Item(39, 5).Value = "  -0.98%  "
Set-TextValue $ws.Cells.Item(40, 4) "1.218.62"
$ws.Cells.Item(40, 5).Value = "  -2.53%  "
Set-TextValue $ws.Cells.Item(41, 4) "6.497"
$ws.Cells.Item(41, 5).Value = "  +5.20%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.9090"
$ws.Cells.Item(42, 5).Value = "  -0.09%  "
$ws.Cells.Item(43, 5).Value = "  +0.02%  "
Set-TextValue $ws.Cells.Item(44, 4) "2.013.51"
$ws.Cells.Item(44, 5).Value = "  +0.02%  "
Set-TextValue $ws.Cells.Item(45, 4) "101.77"
$ws.Cells.Item(45, 5).Value = "  -0.04%  "
Set-TextValue $ws.Cells.Item(46, 4) "66.31"
$ws.Cells.Item(46, 5).Value = "  +0.29%  "
Set-TextValue $ws.Cells.Item(49, 4) "0.4014"
$ws.Cells.Item(49, 5).Value = "  -0.28%  "
Set-TextValue $ws.Cells.Item(50, 4) "8.980"
$ws.Cells.Item(50, 5).Value = "  -1.56%  "
Set-TextValue $ws.Cells.Item(51, 4) "1.680"
$ws.Cells.Item(51, 5).Value = "  -0.29%  "

# Rows 47 and 48 swapped rank order: BabyDogeCoin moved above Aptos
Set-TextValue $ws.Cells.Item(47, 2) "BabyDogeCoin"
Set-TextValue $ws.Cells.Item(47, 3) "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Cells.Item(47, 4) "0.00000000121"
$ws.Cells.Item(47, 5).Value = "  +3.23%  "

Set-TextValue $ws.Cells.Item(48, 2) "Aptos"
Set-TextValue $ws.Cells.Item(48, 3) "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Cells.Item(48, 4) "7.115"
$ws.Cells.Item(48, 5).Value = "  +0.34%  "
